$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Value = "Done"

$ws.Range("E10").Select()
$excel.ActiveWindow.ScrollColumn = 3
